# Applies the upstream edit to PARGT_Manual_Linux.docx:
#   1. The stray "_GoBack" bookmark that used to sit inside the title
#      ("(V|_GoBack|ersion 1.0)") is removed.
#   2. The sentence "... available in the folder- “testing data”." becomes
#      "... available in the folder- “test examples”.", and the "_GoBack"
#      bookmark reappears right after the closing curly quote (i.e. between
#      “test examples” and the trailing period), splitting that sentence
#      into the run layout the diff expects:
#        ["Please note that t"] ["here are ... folder- "] ["“test examples”"]
#        <bookmarkStart/End _GoBack> ["."]

$d = $word.ActiveDocument

# --- Step 1: drop the old "_GoBack" bookmark (wherever it currently lives) ---
if ($d.Bookmarks.Exists("_GoBack")) {
    $d.Bookmarks.Item("_GoBack").Delete()
}

# Helper: force a run split at a character position by dropping a throwaway
# bookmark there and immediately deleting it again - the bookmark's start/end
# boundary is enough to make the engine split the enclosing run in two, and
# that split persists even after the temporary bookmark itself is removed.
function AddSplit($pos) {
    $d.Bookmarks.Add("ZZTEMPSPLIT", $d.Range($pos, $pos)) | Out-Null
    $d.Bookmarks.Item("ZZTEMPSPLIT").Delete()
}

# --- Step 2: locate the sentence and figure out the needed split points ---
# One unambiguous Find across "Please note that t" + "here are ... folder- "
# gives us both boundaries we need in a single shot (the phrase spans the
# existing run break between those two runs).
$sentence = $d.Content
$anchorText = "Please note that there are some example test sequences available in the folder- "
$foundSentence = $sentence.Find.Execute($anchorText, $false, $false, $false, $false, $false, $true, 1, $false, "", 0)
if (-not $foundSentence) {
    throw "Could not locate the target sentence (folder- ...)."
}
$sentenceStart = $sentence.Start
$splitBeforeQuote = $sentence.End
$splitAfterPleaseNote = $sentenceStart + "Please note that t".Length

# --- Step 3: replace the quoted phrase "testing data" with "test examples" ---
$inner = $d.Content
$foundInner = $inner.Find.Execute("testing data", $false, $false, $false, $false, $false, $true, 1, $false, "", 0)
if (-not $foundInner) {
    throw "Could not locate the quoted text 'testing data'."
}
$innerStart = $inner.Start
$inner.Text = "test examples"
$newInnerEnd = $innerStart + ("test examples").Length

# --- Step 4: recreate the run boundaries the diff expects ---
# "Please note that t" | "here are ... folder- "
AddSplit $splitAfterPleaseNote
# "... folder- " | "“test examples”"
AddSplit $splitBeforeQuote

# --- Step 5: re-insert "_GoBack" right after the closing curly quote, ---
#             i.e. between "“test examples”" and the final "."
$closeQuotePos = $newInnerEnd + 1
$d.Bookmarks.Add("_GoBack", $d.Range($closeQuotePos, $closeQuotePos)) | Out-Null

Write-Host "Edit applied successfully."
